$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.804.54"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.34"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.94"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.032"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4411"
$ws.Range("E7").Value = "  +0.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3820"
$ws.Range("E8").Value = "  +1.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07436"
$ws.Range("E9").Value = "  +0.39%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8882"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.64"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.864.15"
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.533"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.748"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07205"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "85.68"
$ws.Range("E16").Value = "  +3.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.040"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009103"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.033"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.57"
$ws.Range("E20").Value = "  +0.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.811.00"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.297"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.29"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.098.28"
$ws.Range("E24").Value = "  +1.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.062"
$ws.Range("E25").Value = "  +6.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.02"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.77"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.008"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.367"
$ws.Range("E29").Value = "  +1.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.67"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09103"
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7766"
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.215"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.023"
$ws.Range("E34").Value = "  +4.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.612"
$ws.Range("E35").Value = "  +2.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.034"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.155"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01987"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05325"
$ws.Range("E39").Value = "  +0.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.871"
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5216"
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.939"
$ws.Range("E42").Value = "  +3.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1677"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.797"
$ws.Range("E44").Value = "  +2.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "110.87"
$ws.Range("E45").Value = "  +1.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.036"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06588"
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.719"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4735"
$ws.Range("E50").Value = "  +1.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.890"
$ws.Range("E51").Value = "  +0.13%  "
